$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 18491262.629693
$ws.Range("D2").Value = 51.753834

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 6332519.424742
$ws.Range("D3").Value = 8.861810999999999
$ws.Range("E3").Value = 0.000198

# Row 4 - Residuals
$ws.Range("B4").Value = 79676253.10051399
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -278.539468
$ws.Range("H5").Value = -539.593939
$ws.Range("I5").Value = -17.484997
$ws.Range("J5").Value = 0.033382

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 80.581307
$ws.Range("H6").Value = -200.519084
$ws.Range("I6").Value = 361.681698
$ws.Range("J6").Value = 0.777482

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 359.120775
$ws.Range("H7").Value = 146.605956
$ws.Range("I7").Value = 571.635594
$ws.Range("J7").Value = 0.000266
